$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.137.38"
$ws.Range("E2").Value = "  -0.87%  "
$ws.Range("D3").Value = "1.558.39"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9971"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9988"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "290.98"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3944"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.99%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3248"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.86%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "43.04"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07355"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.49%  "
$ws.Range("E11").Value = "  -4.08%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9976"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.31%  "
$ws.Range("E13").Value = "  -6.59%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.00001162"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.79%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.664"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.730"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.22%  "
$ws.Range("D17").Value = "1.559.60"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06603"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.61%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "84.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9981"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.364"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.90"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.99%  "
$ws.Range("D24").Value = "22.144.38"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.339"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.488"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "148.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.878"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("D30").Value = "1.735.56"
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "119.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.066"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.748"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08399"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("E35").Value = "  -2.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.627"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -14.40%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06264"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.13%  "
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.201"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2084"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.217"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.11%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "10.83"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9968"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.33%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5877"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.17"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.762"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5634"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.922"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "118.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.95%  "
$ws.Range("E50").Value = "  -2.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06867"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.42%  "
